# Edit the Avancement dashboard workbook:
# - Set all progress values (C5:C8, C10:C14) to 100% (1)
# - Clear the comment cells (D5, D10, D13, D14) that referenced old shared strings
# - Update the active selection to K12
# - The F4 AVERAGE formula will recompute to 1 automatically

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Dashboard")

# Set progress column values to 100% (1) for rows 5-8 and 10-14
$ws.Range("C5").Value = 1
$ws.Range("C6").Value = 1
$ws.Range("C7").Value = 1
$ws.Range("C8").Value = 1
$ws.Range("C10").Value = 1
$ws.Range("C11").Value = 1
$ws.Range("C12").Value = 1
$ws.Range("C13").Value = 1
$ws.Range("C14").Value = 1

# Clear the comment/text cells that are no longer needed
$ws.Range("D5").Value = $null
$ws.Range("D10").Value = $null
$ws.Range("D13").Value = $null
$ws.Range("D14").Value = $null

# Recalculate so the AVERAGE formula in F4 reflects the new values
$excel.Calculate()

# Update the selected cell/active selection to K12
$ws.Activate()
$ws.Range("K12").Select()
